$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.267206788063049
$ws.Range("B1").Value = 1.823354363441467
$ws.Range("C1").Value = 6.810254096984863
$ws.Range("D1").Value = 1.60743772983551
$ws.Range("E1").Value = 0.9494442343711853
